$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: copy format from G1 (bold, bordered, centered) then set value
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Data values for H2:H12 ("Save" column)
$values = @(1, 0, 1, 1, 0, 1, 0, 0, 1, 0, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
